# This document had previously round-tripped through SharePoint, which
# attaches SharePoint-only custom XML data stores to the package:
#   - a content-type / column schema   (ns: .../2006/metadata/contentType)
#   - a document-library form template (ns: .../sharepoint/v3/.../forms)
# "update manual and rebuild" drops that SharePoint metadata again, i.e.
# every non-built-in CustomXMLPart should be removed from the document
# (Word automatically drops the now-unreferenced itemProps part, the
# customXml relationships and the related content-type overrides).

$d = $word.ActiveDocument

$sharePointNamespaces = @(
    "http://schemas.microsoft.com/office/2006/metadata/contentType",
    "http://schemas.microsoft.com/sharepoint/v3/contenttype/forms"
)

# Preferred approach: ask Word directly for the parts in those
# namespaces and delete every match (iterate back-to-front because
# deleting re-indexes the returned collection).
foreach ($ns in $sharePointNamespaces) {
    $matches = $d.CustomXMLParts.SelectByNamespace($ns)
    if ($matches -ne $null) {
        for ($i = $matches.Count; $i -ge 1; $i--) {
            $matches.Item($i).Delete()
        }
    }
}

# Belt-and-braces: walk every CustomXMLPart still attached to the
# document and remove anything that isn't one of Word's own built-in
# parts (covers hosts where SelectByNamespace isn't available).
$parts = $d.CustomXMLParts
for ($i = $parts.Count; $i -ge 1; $i--) {
    $part = $parts.Item($i)
    if (-not $part.BuiltIn) {
        $part.Delete()
    }
}

Write-Output ("CustomXMLParts remaining: " + $d.CustomXMLParts.Count)
